$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44243
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 1200
$ws.Cells.Item(2, 11).Value = 300
$ws.Cells.Item(2, 12).Value = 320
$ws.Cells.Item(2, 13).Value = 310
$ws.Cells.Item(2, 15).Value = "Perú"
$ws.Cells.Item(2, 16).Value = 310

$ws.Cells.Item(3, 4).Value = 44243
$ws.Cells.Item(3, 9).Value = "Segunda"
$ws.Cells.Item(3, 10).Value = 800
$ws.Cells.Item(3, 11).Value = 300
$ws.Cells.Item(3, 12).Value = 320
$ws.Cells.Item(3, 13).Value = 310
$ws.Cells.Item(3, 15).Value = "Perú"
$ws.Cells.Item(3, 16).Value = 310

$ws.Cells.Item(4, 4).Value = 44251
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 1200
$ws.Cells.Item(4, 11).Value = 250
$ws.Cells.Item(4, 12).Value = 280
$ws.Cells.Item(4, 13).Value = 265
$ws.Cells.Item(4, 15).Value = "Perú"
$ws.Cells.Item(4, 16).Value = 265

$ws.Cells.Item(5, 4).Value = 44168
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 1700
$ws.Cells.Item(5, 11).Value = 430
$ws.Cells.Item(5, 12).Value = 450
$ws.Cells.Item(5, 13).Value = 440
$ws.Cells.Item(5, 15).Value = "Perú"
$ws.Cells.Item(5, 16).Value = 440

$ws.Cells.Item(6, 4).Value = 44650
$ws.Cells.Item(6, 9).Value = "Segunda"
$ws.Cells.Item(6, 10).Value = 1000
$ws.Cells.Item(6, 11).Value = 325
$ws.Cells.Item(6, 12).Value = 350
$ws.Cells.Item(6, 13).Value = 338
$ws.Cells.Item(6, 15).Value = "Perú"
$ws.Cells.Item(6, 16).Value = 338

$ws.Cells.Item(7, 4).Value = 44229
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 1200
$ws.Cells.Item(7, 11).Value = 230
$ws.Cells.Item(7, 12).Value = 250
$ws.Cells.Item(7, 13).Value = 240
$ws.Cells.Item(7, 15).Value = "Perú"
$ws.Cells.Item(7, 16).Value = 240

$ws.Cells.Item(8, 4).Value = 44659
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 1000
$ws.Cells.Item(8, 11).Value = 350
$ws.Cells.Item(8, 12).Value = 380
$ws.Cells.Item(8, 13).Value = 365
$ws.Cells.Item(8, 15).Value = "Perú"
$ws.Cells.Item(8, 16).Value = 365

$ws.Cells.Item(9, 4).Value = 44301
$ws.Cells.Item(9, 9).Value = "Segunda"
$ws.Cells.Item(9, 10).Value = 900
$ws.Cells.Item(9, 11).Value = 280
$ws.Cells.Item(9, 12).Value = 300
$ws.Cells.Item(9, 13).Value = 290
$ws.Cells.Item(9, 15).Value = "Perú"
$ws.Cells.Item(9, 16).Value = 290

$ws.Cells.Item(10, 4).Value = 44253
$ws.Cells.Item(10, 9).Value = "Segunda"
$ws.Cells.Item(10, 10).Value = 1200
$ws.Cells.Item(10, 11).Value = 270
$ws.Cells.Item(10, 12).Value = 280
$ws.Cells.Item(10, 13).Value = 275
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 275

$ws.Cells.Item(11, 4).Value = 44580
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 1200
$ws.Cells.Item(11, 11).Value = 380
$ws.Cells.Item(11, 12).Value = 400
$ws.Cells.Item(11, 13).Value = 390
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 390

$ws.Cells.Item(12, 4).Value = 44224
$ws.Cells.Item(12, 9).Value = "Segunda"
$ws.Cells.Item(12, 10).Value = 1200
$ws.Cells.Item(12, 11).Value = 230
$ws.Cells.Item(12, 12).Value = 250
$ws.Cells.Item(12, 13).Value = 240
$ws.Cells.Item(12, 15).Value = "Perú"
$ws.Cells.Item(12, 16).Value = 240

$ws.Cells.Item(13, 4).Value = 44224
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 200
$ws.Cells.Item(13, 12).Value = 230
$ws.Cells.Item(13, 13).Value = 215
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 215

$ws.Cells.Item(14, 4).Value = 44175
$ws.Cells.Item(14, 9).Value = "Segunda"
$ws.Cells.Item(14, 10).Value = 1200
$ws.Cells.Item(14, 11).Value = 400
$ws.Cells.Item(14, 12).Value = 430
$ws.Cells.Item(14, 13).Value = 415
$ws.Cells.Item(14, 15).Value = "Perú"
$ws.Cells.Item(14, 16).Value = 415

$ws.Cells.Item(15, 4).Value = 44231
$ws.Cells.Item(15, 9).Value = "Segunda"
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 180
$ws.Cells.Item(15, 12).Value = 200
$ws.Cells.Item(15, 13).Value = 190
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 190

$ws.Cells.Item(16, 4).Value = 44602
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 1300
$ws.Cells.Item(16, 11).Value = 350
$ws.Cells.Item(16, 12).Value = 380
$ws.Cells.Item(16, 13).Value = 365
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 365

$ws.Cells.Item(17, 4).Value = 44602
$ws.Cells.Item(17, 9).Value = "Segunda"
$ws.Cells.Item(17, 10).Value = 900
$ws.Cells.Item(17, 11).Value = 300
$ws.Cells.Item(17, 12).Value = 330
$ws.Cells.Item(17, 13).Value = 315
$ws.Cells.Item(17, 15).Value = "Perú"
$ws.Cells.Item(17, 16).Value = 315

$ws.Cells.Item(18, 4).Value = 44586
$ws.Cells.Item(18, 9).Value = "Tercera"
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 330
$ws.Cells.Item(18, 12).Value = 350
$ws.Cells.Item(18, 13).Value = 340
$ws.Cells.Item(18, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(18, 16).Value = 340

$ws.Cells.Item(19, 4).Value = 44201
$ws.Cells.Item(19, 9).Value = "Segunda"
$ws.Cells.Item(19, 10).Value = 1800
$ws.Cells.Item(19, 11).Value = 250
$ws.Cells.Item(19, 12).Value = 270
$ws.Cells.Item(19, 13).Value = 260
$ws.Cells.Item(19, 15).Value = "Perú"
$ws.Cells.Item(19, 16).Value = 260

$ws.Cells.Item(20, 4).Value = 44609
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 1200
$ws.Cells.Item(20, 11).Value = 280
$ws.Cells.Item(20, 12).Value = 300
$ws.Cells.Item(20, 13).Value = 290
$ws.Cells.Item(20, 15).Value = "Perú"
$ws.Cells.Item(20, 16).Value = 290

$ws.Cells.Item(21, 4).Value = 44523
$ws.Cells.Item(21, 9).Value = "Segunda"
$ws.Cells.Item(21, 10).Value = 1000
$ws.Cells.Item(21, 11).Value = 550
$ws.Cells.Item(21, 12).Value = 580
$ws.Cells.Item(21, 13).Value = 565
$ws.Cells.Item(21, 15).Value = "Perú"
$ws.Cells.Item(21, 16).Value = 565

$ws.Cells.Item(22, 4).Value = 44176
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 1300
$ws.Cells.Item(22, 11).Value = 350
$ws.Cells.Item(22, 12).Value = 400
$ws.Cells.Item(22, 13).Value = 375
$ws.Cells.Item(22, 15).Value = "Perú"
$ws.Cells.Item(22, 16).Value = 375

$ws.Cells.Item(23, 4).Value = 44172
$ws.Cells.Item(23, 9).Value = "Segunda"
$ws.Cells.Item(23, 10).Value = 1600
$ws.Cells.Item(23, 11).Value = 400
$ws.Cells.Item(23, 12).Value = 420
$ws.Cells.Item(23, 13).Value = 410
$ws.Cells.Item(23, 15).Value = "Perú"
$ws.Cells.Item(23, 16).Value = 410

$ws.Cells.Item(24, 4).Value = 44566
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 1000
$ws.Cells.Item(24, 11).Value = 300
$ws.Cells.Item(24, 12).Value = 320
$ws.Cells.Item(24, 13).Value = 310
$ws.Cells.Item(24, 15).Value = "Perú"
$ws.Cells.Item(24, 16).Value = 310

$ws.Cells.Item(25, 4).Value = 44547
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 1200
$ws.Cells.Item(25, 11).Value = 350
$ws.Cells.Item(25, 12).Value = 370
$ws.Cells.Item(25, 13).Value = 360
$ws.Cells.Item(25, 15).Value = "Perú"
$ws.Cells.Item(25, 16).Value = 360

$ws.Cells.Item(26, 4).Value = 44166
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 1700
$ws.Cells.Item(26, 11).Value = 500
$ws.Cells.Item(26, 12).Value = 530
$ws.Cells.Item(26, 13).Value = 515
$ws.Cells.Item(26, 15).Value = "Perú"
$ws.Cells.Item(26, 16).Value = 515

$ws.Cells.Item(27, 4).Value = 44589
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 900
$ws.Cells.Item(27, 11).Value = 325
$ws.Cells.Item(27, 12).Value = 350
$ws.Cells.Item(27, 13).Value = 338
$ws.Cells.Item(27, 15).Value = "Perú"
$ws.Cells.Item(27, 16).Value = 338

$ws.Cells.Item(28, 4).Value = 44575
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 1200
$ws.Cells.Item(28, 11).Value = 380
$ws.Cells.Item(28, 12).Value = 400
$ws.Cells.Item(28, 13).Value = 390
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value = 390

$ws.Cells.Item(29, 4).Value = 44530
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 1300
$ws.Cells.Item(29, 11).Value = 450
$ws.Cells.Item(29, 12).Value = 480
$ws.Cells.Item(29, 13).Value = 465
$ws.Cells.Item(29, 15).Value = "Perú"
$ws.Cells.Item(29, 16).Value = 465

$ws.Cells.Item(30, 4).Value = 44160
$ws.Cells.Item(30, 9).Value = "Segunda"
$ws.Cells.Item(30, 10).Value = 2000
$ws.Cells.Item(30, 11).Value = 500
$ws.Cells.Item(30, 12).Value = 550
$ws.Cells.Item(30, 13).Value = 525
$ws.Cells.Item(30, 15).Value = "Perú"
$ws.Cells.Item(30, 16).Value = 525

$ws.Cells.Item(31, 4).Value = 44202
$ws.Cells.Item(31, 9).Value = "Segunda"
$ws.Cells.Item(31, 10).Value = 1300
$ws.Cells.Item(31, 11).Value = 230
$ws.Cells.Item(31, 12).Value = 250
$ws.Cells.Item(31, 13).Value = 240
$ws.Cells.Item(31, 15).Value = "Perú"
$ws.Cells.Item(31, 16).Value = 240

$ws.Cells.Item(32, 4).Value = 44162
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 900
$ws.Cells.Item(32, 11).Value = 500
$ws.Cells.Item(32, 12).Value = 550
$ws.Cells.Item(32, 13).Value = 525
$ws.Cells.Item(32, 15).Value = "Perú"
$ws.Cells.Item(32, 16).Value = 525

$ws.Cells.Item(33, 4).Value = 44162
$ws.Cells.Item(33, 9).Value = "Segunda"
$ws.Cells.Item(33, 10).Value = 1200
$ws.Cells.Item(33, 11).Value = 500
$ws.Cells.Item(33, 12).Value = 550
$ws.Cells.Item(33, 13).Value = 525
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 525

$ws.Cells.Item(34, 4).Value = 44214
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 1200
$ws.Cells.Item(34, 11).Value = 400
$ws.Cells.Item(34, 12).Value = 450
$ws.Cells.Item(34, 13).Value = 425
$ws.Cells.Item(34, 15).Value = "Perú"
$ws.Cells.Item(34, 16).Value = 425

$ws.Cells.Item(35, 4).Value = 44603
$ws.Cells.Item(35, 9).Value = "Tercera"
$ws.Cells.Item(35, 10).Value = 300
$ws.Cells.Item(35, 11).Value = 280
$ws.Cells.Item(35, 12).Value = 300
$ws.Cells.Item(35, 13).Value = 290
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 290

$ws.Cells.Item(36, 4).Value = 44217
$ws.Cells.Item(36, 9).Value = "Segunda"
$ws.Cells.Item(36, 10).Value = 1600
$ws.Cells.Item(36, 11).Value = 300
$ws.Cells.Item(36, 12).Value = 350
$ws.Cells.Item(36, 13).Value = 325
$ws.Cells.Item(36, 15).Value = "Perú"
$ws.Cells.Item(36, 16).Value = 325
